# Generate Report for handoff
# The "6bf2a2b1-212a-4798-bf92-426f6248ecf6.md" file has finished its
# handoff pass and moved into translation, while the files already queued
# for handoff (b346fae3, bb845ad7, 7cd9e49e, 8ba81ca8, a24df2d9) get a
# fresh "Latest Handoff Datetime" stamp from this handoff run.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: row for 6bf2a2b1-212a-4798-bf92-426f6248ecf6.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = $newStatus
$zhcn.Range("D2").Value = "2016-01-28 09:57:27"
$zhcn.Range("D3").Value = "2016-01-28 09:57:27"
$zhcn.Range("D5").Value = "2016-01-28 09:57:27"
$zhcn.Range("D6").Value = "2016-01-28 09:57:27"
$zhcn.Range("D7").Value = "2016-01-28 09:57:27"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = $newStatus
$dede.Range("D2").Value = "2016-01-28 09:57:41"
$dede.Range("D3").Value = "2016-01-28 09:57:41"
$dede.Range("D5").Value = "2016-01-28 09:57:41"
$dede.Range("D6").Value = "2016-01-28 09:57:41"
$dede.Range("D7").Value = "2016-01-28 09:57:41"
